$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row and title-case municipality/state names
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B4').Value = 'Rincón De Romos'
$ws.Range('B18').Value = 'Amatenango De La Frontera'
$ws.Range('B21').Value = 'Bejucal De Ocampo'
$ws.Range('B30').Value = 'Comitán De Domínguez'
$ws.Range('B44').Value = 'Marqués De Comillas'
$ws.Range('B45').Value = 'Mazapa De Madero'
$ws.Range('B49').Value = 'Salto De Agua'
$ws.Range('A81').Value = 'Ciudad De México'
$ws.Range('B85').Value = 'Cuajimalpa De Morelos'
$ws.Range('B102').Value = 'San Juan Del Río'
$ws.Range('A105').Value = 'Estado De México'
$ws.Range('B105').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B107').Value = 'Almoloya De Alquisiras'
$ws.Range('B108').Value = 'Almoloya De Juárez'
$ws.Range('B109').Value = 'Almoloya Del Río'
$ws.Range('B114').Value = 'Atizapán De Zaragoza'
$ws.Range('B120').Value = 'Chapa De Mota'
$ws.Range('B127').Value = 'Ecatepec De Morelos'
$ws.Range('B131').Value = 'Ixtapan De La Sal'
$ws.Range('B132').Value = 'Ixtapan Del Oro'
$ws.Range('B140').Value = 'Naucalpan De Juárez'
$ws.Range('B148').Value = 'San Felipe Del Progreso'
$ws.Range('B149').Value = 'San Simón De Guerrero'
$ws.Range('B158').Value = 'Tenango Del Valle'
$ws.Range('B165').Value = 'Tlalnepantla De Baz'
$ws.Range('B168').Value = 'Valle De Bravo'
$ws.Range('B169').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B170').Value = 'Villa De Allende'
$ws.Range('B181').Value = 'San Miguel De Allende'
$ws.Range('B182').Value = 'Apaseo El Alto'
$ws.Range('B183').Value = 'Apaseo El Grande'
$ws.Range('B189').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B193').Value = 'Jaral Del Progreso'
$ws.Range('B201').Value = 'San Diego De La Unión'
$ws.Range('B203').Value = 'San Francisco Del Rincón'
$ws.Range('B204').Value = 'San Luis De La Paz'
$ws.Range('B205').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B207').Value = 'Silao De La Victoria'
$ws.Range('B211').Value = 'Valle De Santiago'
$ws.Range('B217').Value = 'Acapulco De Juárez'
$ws.Range('B218').Value = 'Alcozauca De Guerrero'
$ws.Range('B222').Value = 'Atlamajalcingo Del Monte'
$ws.Range('B224').Value = 'Atoyac De Álvarez'
$ws.Range('B225').Value = 'Ayutla De Los Libres'
$ws.Range('B227').Value = 'Chilapa De Álvarez'
$ws.Range('B228').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B229').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B232').Value = 'Coyuca De Benítez'
$ws.Range('B233').Value = 'Coyuca De Catalán'
$ws.Range('B237').Value = 'Cuetzala Del Progreso'
$ws.Range('B238').Value = 'Cutzamala De Pinzón'
$ws.Range('B243').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B244').Value = 'Iguala De La Independencia'
$ws.Range('B245').Value = 'Zihuatanejo De Azueta'
$ws.Range('B259').Value = 'Taxco De Alarcón'
$ws.Range('B261').Value = 'Técpan De Galeana'
$ws.Range('B263').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B266').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B267').Value = 'Tlapa De Comonfort'
$ws.Range('B279').Value = 'Atotonilco De Tula'
$ws.Range('B280').Value = 'Atotonilco El Grande'
$ws.Range('B283').Value = 'Cuautepec De Hinojosa'
$ws.Range('B286').Value = 'Huasca De Ocampo'
$ws.Range('B292').Value = 'Mineral De La Reforma'
$ws.Range('B293').Value = 'Mineral Del Chico'
$ws.Range('B295').Value = 'Omitlán De Juárez'
$ws.Range('B296').Value = 'Pachuca De Soto'
$ws.Range('B298').Value = 'Santiago De Anaya'
$ws.Range('B301').Value = 'Tepehuacán De Guerrero'
$ws.Range('B302').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B303').Value = 'Tezontepec De Aldama'
$ws.Range('B306').Value = 'Tula De Allende'
$ws.Range('B307').Value = 'Tulancingo De Bravo'
$ws.Range('B309').Value = 'Zacualtipán De Ángeles'
$ws.Range('B310').Value = 'Zapotlán De Juárez'
$ws.Range('B313').Value = 'Acatlán De Juárez'
$ws.Range('B314').Value = 'Ahualulco De Mercado'
$ws.Range('B315').Value = 'Atotonilco El Alto'
$ws.Range('B316').Value = 'Autlán De Navarro'
$ws.Range('B324').Value = 'Encarnación De Díaz'
$ws.Range('B331').Value = 'La Manzanilla De La Paz'
$ws.Range('B332').Value = 'Lagos De Moreno'
$ws.Range('B336').Value = 'San Juan De Los Lagos'
$ws.Range('B337').Value = 'San Martín De Bolaños'
$ws.Range('B339').Value = 'Tamazula De Gordiano'
$ws.Range('B342').Value = 'Tepatitlán De Morelos'
$ws.Range('B343').Value = 'Tizapán El Alto'
$ws.Range('B398').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B425').Value = 'Tetela Del Volcán'
$ws.Range('B426').Value = 'Tlaltizapán De Zapata'
$ws.Range('B430').Value = 'Zacualpan De Amilpas'
$ws.Range('B434').Value = 'Santa María Del Oro'
$ws.Range('B443').Value = 'San Nicolás De Los Garza'
$ws.Range('B445').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B451').Value = 'Fresnillo De Trujano'
$ws.Range('B452').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B453').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B454').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B455').Value = 'Huautla De Jiménez'
$ws.Range('B456').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B459').Value = 'Mariscala De Juárez'
$ws.Range('B460').Value = 'Mártires De Tacubaya'
$ws.Range('B462').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B465').Value = 'Oaxaca De Juárez'
$ws.Range('B466').Value = 'Ocotlán De Morelos'
$ws.Range('B467').Value = 'Putla Villa De Guerrero'
$ws.Range('B482').Value = 'San Dionisio Del Mar'
$ws.Range('B486').Value = 'San José Del Progreso'
$ws.Range('B488').Value = 'San Juan Bautista Lo De Soto'
$ws.Range('B504').Value = 'San Miguel Del Puerto'
$ws.Range('B506').Value = 'San Pablo Villa De Mitla'
$ws.Range('B513').Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range('B514').Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Range('B528').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B554').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B555').Value = 'Tlacolula De Matamoros'
$ws.Range('B556').Value = 'Tlalixtac De Cabrera'
$ws.Range('B559').Value = 'Villa De Chilapa De Díaz'
$ws.Range('B560').Value = 'Villa De Etla'
$ws.Range('B561').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B562').Value = 'Villa Sola De Vega'
$ws.Range('B564').Value = 'Zimatlán De Álvarez'
$ws.Range('B584').Value = 'Chalchicomula De Sesma'
$ws.Range('B596').Value = 'Cuayuca De Andrade'
$ws.Range('B597').Value = 'Cuetzalan Del Progreso'
$ws.Range('B608').Value = 'Huehuetlán El Chico'
$ws.Range('B609').Value = 'Huehuetlán El Grande'
$ws.Range('B612').Value = 'Huitzilan De Serdán'
$ws.Range('B614').Value = 'Ixcamilpa De Guerrero'
$ws.Range('B616').Value = 'Izúcar De Matamoros'
$ws.Range('B621').Value = 'Los Reyes De Juárez'
$ws.Range('B626').Value = 'Palmar De Bravo'
$ws.Range('B640').Value = 'San Nicolás De Los Ranchos'
$ws.Range('B643').Value = 'San Salvador El Seco'
$ws.Range('B644').Value = 'San Salvador El Verde'
$ws.Range('B651').Value = 'Tecali De Herrera'
$ws.Range('B657').Value = 'Tepanco De López'
$ws.Range('B658').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B662').Value = 'Tepexi De Rodríguez'
$ws.Range('B664').Value = 'Tetela De Ocampo'
$ws.Range('B665').Value = 'Teteles De Avila Castillo'
$ws.Range('B669').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B680').Value = 'Xayacatlán De Bravo'
$ws.Range('B684').Value = 'Xochitlán De Vicente Suárez'
$ws.Range('B694').Value = 'Amealco De Bonfil'
$ws.Range('B695').Value = 'Cadereyta De Montes'
$ws.Range('B698').Value = 'Jalpan De Serra'
$ws.Range('B699').Value = 'Pinal De Amoles'
$ws.Range('B702').Value = 'San Juan Del Río'
$ws.Range('B708').Value = 'Armadillo De Los Infante'
$ws.Range('B709').Value = 'Axtla De Terrazas'
$ws.Range('B710').Value = 'Ciudad Del Maíz'
$ws.Range('B723').Value = 'Villa De Guadalupe'
$ws.Range('B764').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B765').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B767').Value = 'Papalotla De Xicohténcatl'
$ws.Range('B770').Value = 'San Pablo Del Monte'
$ws.Range('B773').Value = 'Tepetitla De Lardizábal'
$ws.Range('B775').Value = 'Tetla De La Solidaridad'
$ws.Range('B785').Value = 'Amatlán De Los Reyes'
$ws.Range('B788').Value = 'Boca Del Río'
$ws.Range('B791').Value = 'Chinampa De Gorostiza'
$ws.Range('B798').Value = 'Cosamaloapan De Carpio'
$ws.Range('B807').Value = 'Ignacio De La Llave'
$ws.Range('B809').Value = 'Ixhuatlán De Madero'
$ws.Range('B810').Value = 'Ixhuatlán Del Café'
$ws.Range('B816').Value = 'Lerdo De Tejada'
$ws.Range('B819').Value = 'Martínez De La Torre'
$ws.Range('B820').Value = 'Medellín De Bravo'
$ws.Range('B824').Value = 'Mixtla De Altamirano'
$ws.Range('B831').Value = 'Paso De Ovejas'
$ws.Range('B832').Value = 'Paso Del Macho'
$ws.Range('B835').Value = 'Poza Rica De Hidalgo'
$ws.Range('B839').Value = 'Soledad De Doblado'
$ws.Range('B843').Value = 'Tatahuicapan De Juárez'
$ws.Range('B857').Value = 'Vega De Alatorre'
$ws.Range('B876').Value = 'Villa De Cos'

# Remove trailing footer/metadata rows (882-886)
$ws.Range("A882:A886").EntireRow.Delete()
